$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '51.606.93'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.782.11'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '352.69'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').Value = '109.00'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -2.49%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +2.26%  '
$ws.Range('D10').Value = '39.66'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('D13').Value = '20.08'
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').Value = '3.219.42'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '2.779.76'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '51.646.65'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').Value = '7.71'
$ws.Range('E19').Value = '  +3.85%  '
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').Value = '69.88'
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').Value = '267.36'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').Value = '2.72'
$ws.Range('D26').Value = '26.11'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = '0.164'
$ws.Range('E28').Value = '  +12.66%  '
$ws.Range('D29').Value = '10.24'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = '37.00'
$ws.Range('E30').Value = '  +7.27%  '
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('E32').Value = '  +7.20%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '5.68'
$ws.Range('E34').Value = '  +7.85%  '
$ws.Range('D35').Value = '0.0453'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D36').Value = '0.0833'
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '18.49'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('E39').Value = '  -3.02%  '
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D43').Value = '22.12'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').Value = '120.18'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('D46').Value = '2.126.94'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('E48').Value = '  +4.23%  '
$ws.Range('D49').Value = '5.43'
$ws.Range('E49').Value = '  -5.06%  '
$ws.Range('D50').Value = '0.905'
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('E51').Value = '  +8.41%  '
